$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into $addr while preserving whatever style that cell
# already carries. Does this by staging the value in a scratch cell (far off
# to the side), copying it, and PasteSpecial-ing *values only* into the
# target so the destination's existing number format / alignment survives.
# ---------------------------------------------------------------------------
$scratch = "AZ1"
function Set-ValueKeepStyle($ws, $addr, $value) {
    $ws.Range($scratch).Value = $value
    $ws.Range($scratch).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# Helper: stamp $addr with the same formatting as $srcAddr (without touching
# the value), then write $value into it while preserving that stamped style.
# Used for brand-new cells that need to inherit a format that doesn't exist
# on them yet (e.g. the 0.000 number style or the right-aligned record style).
# ---------------------------------------------------------------------------
function Set-NewStyledValue($ws, $addr, $srcAddr, $value) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    Set-ValueKeepStyle $ws $addr $value
}

# ===========================================================================
# 1) "Best of all Time (Playoff Era)" table, columns G:J
#    Add the 2020 Alabama game into the first open slot (row 12, which
#    already carries the table's number/record styles on I12/J12), then
#    re-sort the table (rows 3-12) descending by score so it lands in order.
# ===========================================================================
$ws.Range("G12").Value = 2020
$ws.Range("H12").Value = "Alabama"
Set-ValueKeepStyle $ws "I12" 40.866
Set-ValueKeepStyle $ws "J12" "13-0"

$bcsRange = $ws.Range("G2:J12")
$bcsKey = $ws.Range("I2")
$bcsRange.Sort($bcsKey, 2, $null, $null, 1, $null, $null, 1)

# ===========================================================================
# 2) "Worst of all Time (Playoff Era)" table, columns S:V
#    Add the three new 2020 games into fresh rows 10-12 (none of these
#    existed before, so stamp U/V with the table's existing styles first),
#    then re-sort the table (rows 3-12) ascending by score.
# ===========================================================================
$ws.Range("S10").Value = 2020
$ws.Range("T10").Value = "Vanderbilt"
Set-NewStyledValue $ws "U10" "U9" 15.04
Set-NewStyledValue $ws "V10" "V9" "0-9"

$ws.Range("S11").Value = 2020
$ws.Range("T11").Value = "Bowling Green"
Set-NewStyledValue $ws "U11" "U9" 15.61
Set-NewStyledValue $ws "V11" "V9" "0-5"

$ws.Range("S12").Value = 2020
$ws.Range("T12").Value = "FIU"
Set-NewStyledValue $ws "U12" "U9" 15.79
Set-NewStyledValue $ws "V12" "V9" "0-5"

$woatRange = $ws.Range("S2:V12")
$woatKey = $ws.Range("U2")
$woatRange.Sort($woatKey, 1, $null, $null, 1, $null, $null, 1)

# ===========================================================================
# 3) Running averages now need to cover the extra row each table picked up.
# ===========================================================================
$ws.Range("K1").Formula = "=AVERAGE(I3:I12)"
$ws.Range("W1").Formula = "=AVERAGE(U3:U12)"

# ===========================================================================
# 4) Clean up the scratch cell, fix selection + the filter-database name.
# ===========================================================================
$ws.Range($scratch).Clear()
$ws.Range("R4").Select()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$S`$2:`$V`$2"
    }
}

Write-Output "done"
